# Add a third row to the "data_info(behavior)" sheet that carries the
# French description / enum labels for each existing column (Operator,
# SampleID, Date, LaboratoryOperatingMode, CriticalApparatusCriticalSoftware,
# CriticalProduct, RawDataPathway). Row 1 holds the column names and row 2
# holds the data-type tag (#string/#date); row 3 adds the new
# description/enum tag for each column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
